$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026327809050455
$ws.Cells.Item(2, 4).Value = 1.033934170463273
$ws.Cells.Item(2, 5).Value = 1.026562891491227
$ws.Cells.Item(2, 6).Value = 1.041029263270844
$ws.Cells.Item(2, 9).Value = 1.030972423753723
$ws.Cells.Item(2, 10).Value = 1.031491865453045
$ws.Cells.Item(2, 11).Value = 1.036735092774435
$ws.Cells.Item(2, 12).Value = 1.02938517384718
$ws.Cells.Item(2, 13).Value = 1.043809936069994

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027284868655318
$ws.Cells.Item(3, 4).Value = 1.034652058386466
$ws.Cells.Item(3, 5).Value = 1.027375804996841
$ws.Cells.Item(3, 6).Value = 1.041943570300468
$ws.Cells.Item(3, 9).Value = 1.031104258443841
$ws.Cells.Item(3, 10).Value = 1.032088599008431
$ws.Cells.Item(3, 11).Value = 1.037262307209987
$ws.Cells.Item(3, 12).Value = 1.030005613419861
$ws.Cells.Item(3, 13).Value = 1.044534514629731

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.02790469430447
$ws.Cells.Item(4, 4).Value = 1.035116844112886
$ws.Cells.Item(4, 5).Value = 1.027902660592929
$ws.Cells.Item(4, 6).Value = 1.042535956090852
$ws.Cells.Item(4, 9).Value = 1.031188277057087
$ws.Cells.Item(4, 10).Value = 1.032474676121379
$ws.Cells.Item(4, 11).Value = 1.037603044275085
$ws.Cells.Item(4, 12).Value = 1.030407276591723
$ws.Cells.Item(4, 13).Value = 1.045003503623353

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028165398288454
$ws.Cells.Item(5, 4).Value = 1.035312301825979
$ws.Cells.Item(5, 5).Value = 1.028124351883175
$ws.Cells.Item(5, 6).Value = 1.042785177462964
$ws.Cells.Item(5, 9).Value = 1.031223289755673
$ws.Cells.Item(5, 10).Value = 1.032636970101201
$ws.Cells.Item(5, 11).Value = 1.037746191859674
$ws.Cells.Item(5, 12).Value = 1.030576181941058
$ws.Cells.Item(5, 13).Value = 1.045200698402348

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028209179187396
$ws.Cells.Item(6, 4).Value = 1.035345123617064
$ws.Cells.Item(6, 5).Value = 1.028161586584383
$ws.Cells.Item(6, 6).Value = 1.04282703348028
$ws.Cells.Item(6, 9).Value = 1.031229150421788
$ws.Cells.Item(6, 10).Value = 1.032664219192287
$ws.Cells.Item(6, 11).Value = 1.03777022117406
$ws.Cells.Item(6, 12).Value = 1.030604544567486
$ws.Cells.Item(6, 13).Value = 1.045233810101291

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027908177337089
$ws.Cells.Item(7, 4).Value = 1.035119455587451
$ws.Cells.Item(7, 5).Value = 1.027905622053969
$ws.Cells.Item(7, 6).Value = 1.042539285484224
$ws.Cells.Item(7, 9).Value = 1.031188746112258
$ws.Cells.Item(7, 10).Value = 1.032476844753196
$ws.Cells.Item(7, 11).Value = 1.037604957407659
$ws.Cells.Item(7, 12).Value = 1.030409533333274
$ws.Cells.Item(7, 13).Value = 1.045006138426275

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.026651138580373
$ws.Cells.Item(8, 4).Value = 1.034176728455599
$ws.Cells.Item(8, 5).Value = 1.026837443372791
$ws.Cells.Item(8, 6).Value = 1.041338097704046
$ws.Cells.Item(8, 9).Value = 1.031017244102342
$ws.Cells.Item(8, 10).Value = 1.031693543971738
$ws.Cells.Item(8, 11).Value = 1.03691335075156
$ws.Cells.Item(8, 12).Value = 1.029594812809339
$ws.Cells.Item(8, 13).Value = 1.044054781396932

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02444028580838
$ws.Cells.Item(9, 4).Value = 1.032517611180028
$ws.Cells.Item(9, 5).Value = 1.024961718060623
$ws.Cells.Item(9, 6).Value = 1.039227398532732
$ws.Cells.Item(9, 9).Value = 1.030705201477232
$ws.Cells.Item(9, 10).Value = 1.03031293423995
$ws.Cells.Item(9, 11).Value = 1.035691590962966
$ws.Cells.Item(9, 12).Value = 1.028160729648756
$ws.Cells.Item(9, 13).Value = 1.042379483940625

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.022969271401684
$ws.Cells.Item(10, 4).Value = 1.031413023755124
$ws.Cells.Item(10, 5).Value = 1.023715714828007
$ws.Cells.Item(10, 6).Value = 1.037824351427936
$ws.Cells.Item(10, 9).Value = 1.030490594060963
$ws.Cells.Item(10, 10).Value = 1.02939236299889
$ws.Cells.Item(10, 11).Value = 1.034875091508912
$ws.Cells.Item(10, 12).Value = 1.027205785615194
$ws.Cells.Item(10, 13).Value = 1.041263444447964

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022333002344295
$ws.Cells.Item(11, 4).Value = 1.030935096906643
$ws.Cells.Item(11, 5).Value = 1.023177261423274
$ws.Cells.Item(11, 6).Value = 1.037217803750861
$ws.Cells.Item(11, 9).Value = 1.030396113085392
$ws.Cells.Item(11, 10).Value = 1.028993718768648
$ws.Cells.Item(11, 11).Value = 1.034521079020396
$ws.Cells.Item(11, 12).Value = 1.026792561138263
$ws.Cells.Item(11, 13).Value = 1.040780398568353

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022096768213192
$ws.Cells.Item(12, 4).Value = 1.0307576301433
$ws.Cells.Item(12, 5).Value = 1.022977418603912
$ws.Cells.Item(12, 6).Value = 1.036992653754951
$ws.Cells.Item(12, 9).Value = 1.030360785616719
$ws.Cells.Item(12, 10).Value = 1.028845640889971
$ws.Cells.Item(12, 11).Value = 1.03438951454234
$ws.Cells.Item(12, 12).Value = 1.026639113164868
$ws.Cells.Item(12, 13).Value = 1.040601006099229

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022147436492011
$ws.Cells.Item(13, 4).Value = 1.030795694784428
$ws.Cells.Item(13, 5).Value = 1.023020278180429
$ws.Cells.Item(13, 6).Value = 1.037040942435956
$ws.Cells.Item(13, 9).Value = 1.03036837401521
$ws.Cells.Item(13, 10).Value = 1.028877404247788
$ws.Cells.Item(13, 11).Value = 1.034417738651237
$ws.Cells.Item(13, 12).Value = 1.026672026358896
$ws.Cells.Item(13, 13).Value = 1.040639484912228

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022313473019909
$ws.Cells.Item(14, 4).Value = 1.030920426286735
$ws.Cells.Item(14, 5).Value = 1.023160739025567
$ws.Cells.Item(14, 6).Value = 1.037199189737808
$ws.Cells.Item(14, 9).Value = 1.030393197656459
$ws.Cells.Item(14, 10).Value = 1.028981478674111
$ws.Cells.Item(14, 11).Value = 1.034510205250393
$ws.Cells.Item(14, 12).Value = 1.026779876218174
$ws.Cells.Item(14, 13).Value = 1.040765569265248

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022415787481829
$ws.Cells.Item(15, 4).Value = 1.030997285018189
$ws.Cells.Item(15, 5).Value = 1.023247303194682
$ws.Cells.Item(15, 6).Value = 1.037296710886118
$ws.Cells.Item(15, 9).Value = 1.030408461458443
$ws.Cells.Item(15, 10).Value = 1.029045601906751
$ws.Cells.Item(15, 11).Value = 1.034567167931664
$ws.Cells.Item(15, 12).Value = 1.026846331672059
$ws.Cells.Item(15, 13).Value = 1.040843258310763

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023011513093559
$ws.Cells.Item(16, 4).Value = 1.031444750037035
$ws.Cells.Item(16, 5).Value = 1.023751472955808
$ws.Cells.Item(16, 6).Value = 1.037864626790366
$ws.Cells.Item(16, 9).Value = 1.030496831744453
$ws.Cells.Item(16, 10).Value = 1.029418819136944
$ws.Cells.Item(16, 11).Value = 1.034898576490777
$ws.Cells.Item(16, 12).Value = 1.027233215794957
$ws.Cells.Item(16, 13).Value = 1.041295507070554

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023385381352319
$ws.Cells.Item(17, 4).Value = 1.031725532240765
$ws.Cells.Item(17, 5).Value = 1.024068013886138
$ws.Cells.Item(17, 6).Value = 1.038221129082794
$ws.Cells.Item(17, 9).Value = 1.030551848264283
$ws.Cells.Item(17, 10).Value = 1.029652920879071
$ws.Cells.Item(17, 11).Value = 1.035106337156902
$ws.Cells.Item(17, 12).Value = 1.027475971700072
$ws.Cells.Item(17, 13).Value = 1.041579246953994

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.023603518887262
$ws.Cells.Item(18, 4).Value = 1.031889343043635
$ws.Cells.Item(18, 5).Value = 1.024252750410497
$ws.Cells.Item(18, 6).Value = 1.038429165392149
$ws.Cells.Item(18, 9).Value = 1.03058378844283
$ws.Cells.Item(18, 10).Value = 1.029789465427152
$ws.Cells.Item(18, 11).Value = 1.035227475710364
$ws.Cells.Item(18, 12).Value = 1.02761759334359
$ws.Cells.Item(18, 13).Value = 1.041744767476727

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023677909403233
$ws.Cells.Item(19, 4).Value = 1.031945204236763
$ws.Cells.Item(19, 5).Value = 1.024315758314632
$ws.Cells.Item(19, 6).Value = 1.038500116391746
$ws.Cells.Item(19, 9).Value = 1.030594653755203
$ws.Cells.Item(19, 10).Value = 1.029836023044831
$ws.Cells.Item(19, 11).Value = 1.03526877320774
$ws.Cells.Item(19, 12).Value = 1.027665887081754
$ws.Cells.Item(19, 13).Value = 1.041801209026574

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023345261941102
$ws.Cells.Item(20, 4).Value = 1.031695403306565
$ws.Cells.Item(20, 5).Value = 1.024034041318989
$ws.Cells.Item(20, 6).Value = 1.038182869970186
$ws.Cells.Item(20, 9).Value = 1.030545961026537
$ws.Cells.Item(20, 10).Value = 1.02962780428483
$ws.Cells.Item(20, 11).Value = 1.035084051022424
$ws.Cells.Item(20, 12).Value = 1.027449923568114
$ws.Cells.Item(20, 13).Value = 1.041548802287706

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022264576508261
$ws.Cells.Item(21, 4).Value = 1.030883694395592
$ws.Cells.Item(21, 5).Value = 1.02311937230395
$ws.Cells.Item(21, 6).Value = 1.037152585731974
$ws.Cells.Item(21, 9).Value = 1.030385894141298
$ws.Cells.Item(21, 10).Value = 1.028950831442864
$ws.Cells.Item(21, 11).Value = 1.034482978023419
$ws.Cells.Item(21, 12).Value = 1.02674811595503
$ws.Cells.Item(21, 13).Value = 1.040728439664846

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021585711882599
$ws.Cells.Item(22, 4).Value = 1.030373668667639
$ws.Cells.Item(22, 5).Value = 1.022545225884821
$ws.Cells.Item(22, 6).Value = 1.036505666866172
$ws.Cells.Item(22, 9).Value = 1.030283905772787
$ws.Cells.Item(22, 10).Value = 1.028525171040694
$ws.Cells.Item(22, 11).Value = 1.034104663907472
$ws.Cells.Item(22, 12).Value = 1.026307105147278
$ws.Cells.Item(22, 13).Value = 1.040212832490367

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.021945533204537
$ws.Cells.Item(23, 4).Value = 1.030644011395868
$ws.Cells.Item(23, 5).Value = 1.022849501991879
$ws.Cells.Item(23, 6).Value = 1.036848528635186
$ws.Cells.Item(23, 9).Value = 1.030338099332293
$ws.Cells.Item(23, 10).Value = 1.028750823313415
$ws.Cells.Item(23, 11).Value = 1.034305252529201
$ws.Cells.Item(23, 12).Value = 1.026540869892108
$ws.Cells.Item(23, 13).Value = 1.040486147434354

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023363389970664
$ws.Cells.Item(24, 4).Value = 1.031709017166144
$ws.Cells.Item(24, 5).Value = 1.024049391739925
$ws.Cells.Item(24, 6).Value = 1.038200157324313
$ws.Cells.Item(24, 9).Value = 1.030548621679468
$ws.Cells.Item(24, 10).Value = 1.029639153402105
$ws.Cells.Item(24, 11).Value = 1.035094121305599
$ws.Cells.Item(24, 12).Value = 1.027461693516766
$ws.Cells.Item(24, 13).Value = 1.041562558860711

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025011339488399
$ws.Cells.Item(25, 4).Value = 1.032946276347523
$ws.Cells.Item(25, 5).Value = 1.025445854320967
$ws.Cells.Item(25, 6).Value = 1.039772351116418
$ws.Cells.Item(25, 9).Value = 1.030787033984337
$ws.Cells.Item(25, 10).Value = 1.030669888353762
$ws.Cells.Item(25, 11).Value = 1.036007801163984
$ws.Cells.Item(25, 12).Value = 1.028531283282608
$ws.Cells.Item(25, 13).Value = 1.042812448537928
